$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data of rows 3, 4 and 5 (excluding the first few
# identifying columns A:C which stay identical across the three rows):
#   new row3 <- old row4
#   new row4 <- old row5
#   new row5 <- old row3
# Columns affected: D (Fecha) and L:T (Calidad ... Kg/unidad)

$oldRow3 = @{
    D = $ws.Range("D3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    N = $ws.Range("N3").Value2
    O = $ws.Range("O3").Value2
    P = $ws.Range("P3").Value2
    Q = $ws.Range("Q3").Value2
    R = $ws.Range("R3").Value2
    S = $ws.Range("S3").Value2
    T = $ws.Range("T3").Value2
}

$oldRow4 = @{
    D = $ws.Range("D4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    N = $ws.Range("N4").Value2
    O = $ws.Range("O4").Value2
    P = $ws.Range("P4").Value2
    Q = $ws.Range("Q4").Value2
    R = $ws.Range("R4").Value2
    S = $ws.Range("S4").Value2
    T = $ws.Range("T4").Value2
}

$oldRow5 = @{
    D = $ws.Range("D5").Value2
    L = $ws.Range("L5").Value2
    M = $ws.Range("M5").Value2
    N = $ws.Range("N5").Value2
    O = $ws.Range("O5").Value2
    P = $ws.Range("P5").Value2
    Q = $ws.Range("Q5").Value2
    R = $ws.Range("R5").Value2
    S = $ws.Range("S5").Value2
    T = $ws.Range("T5").Value2
}

function Set-RowValues($rowNum, $data) {
    $ws.Range("D$rowNum").Value2 = $data.D
    $ws.Range("L$rowNum").Value2 = $data.L
    $ws.Range("M$rowNum").Value2 = $data.M
    $ws.Range("N$rowNum").Value2 = $data.N
    $ws.Range("O$rowNum").Value2 = $data.O
    $ws.Range("P$rowNum").Value2 = $data.P
    $ws.Range("Q$rowNum").Value2 = $data.Q
    $ws.Range("R$rowNum").Value2 = $data.R
    $ws.Range("S$rowNum").Value2 = $data.S
    $ws.Range("T$rowNum").Value2 = $data.T
}

Set-RowValues 3 $oldRow4
Set-RowValues 4 $oldRow5
Set-RowValues 5 $oldRow3
